$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.837.73"
$ws.Range("E2").Value = "  -5.20%  "

# Row 3
$ws.Range("D3").Value = "3.368.51"
$ws.Range("E3").Value = "  -6.47%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").Value = "563.51"
$ws.Range("E5").Value = "  -5.58%  "

# Row 6
$ws.Range("D6").Value = "184.71"
$ws.Range("E6").Value = "  -8.60%  "

# Row 7
$ws.Range("D7").Value = "0.596"
$ws.Range("E7").Value = "  -5.03%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("D9").Value = "3.363.13"
$ws.Range("E9").Value = "  -6.29%  "

# Row 10
$ws.Range("E10").Value = "  -12.51%  "

# Row 11
$ws.Range("E11").Value = "  -7.25%  "

# Row 12
$ws.Range("D12").Value = "47.90"
$ws.Range("E12").Value = "  -11.02%  "

# Row 13
$ws.Range("D13").Value = "0.0000268"
$ws.Range("E13").Value = "  -11.35%  "

# Row 14
$ws.Range("D14").Value = "8.75"
$ws.Range("E14").Value = "  -9.14%  "

# Row 15
$ws.Range("D15").Value = "3.901.18"

# Row 16
$ws.Range("D16").Value = "606.00"
$ws.Range("E16").Value = "  -10.62%  "

# Row 17
$ws.Range("D17").Value = "66.771.96"
$ws.Range("E17").Value = "  -5.40%  "

# Row 18
$ws.Range("E18").Value = "  -3.60%  "

# Row 19
$ws.Range("D19").Value = "3.366.84"
$ws.Range("E19").Value = "  -6.21%  "

# Row 20
$ws.Range("D20").Value = "17.67"
$ws.Range("E20").Value = "  -7.65%  "

# Row 21
$ws.Range("D21").Value = "11.68"
$ws.Range("E21").Value = "  -8.50%  "

# Row 22
$ws.Range("D22").Value = "0.914"
$ws.Range("E22").Value = "  -8.57%  "

# Row 23
$ws.Range("D23").Value = "17.12"
$ws.Range("E23").Value = "  -8.31%  "

# Row 24
$ws.Range("D24").Value = "5.09"
$ws.Range("E24").Value = "  -3.55%  "

# Row 25
$ws.Range("D25").Value = "95.46"
$ws.Range("E25").Value = "  -13.57%  "

# Row 26
$ws.Range("E26").Value = "  -9.78%  "

# Row 27
$ws.Range("E27").Value = "  -9.45%  "

# Row 28
$ws.Range("E28").Value = "  -9.19%  "

# Row 29
$ws.Range("D29").Value = "8.81"
$ws.Range("E29").Value = "  -12.76%  "

# Row 30
$ws.Range("D30").Value = "30.82"
$ws.Range("E30").Value = "  -9.86%  "

# Row 31
$ws.Range("D31").Value = "6.38"
$ws.Range("E31").Value = "  -10.95%  "

# Row 32
$ws.Range("D32").Value = "3.87"
$ws.Range("E32").Value = "  -13.69%  "

# Row 33
$ws.Range("D33").Value = "11.28"
$ws.Range("E33").Value = "  -8.51%  "

# Row 34
$ws.Range("E34").Value = "  -7.86%  "

# Row 35
$ws.Range("D35").Value = "58.67"
$ws.Range("E35").Value = "  -7.72%  "

# Row 36
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "3.767.28"
$ws.Range("E36").Value = "  -2.72%  "

# Row 37
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "532.42"
$ws.Range("E37").Value = "  +3.95%  "

# Row 38
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.06%  "

# Row 39
$ws.Range("D39").Value = "3.80"
$ws.Range("E39").Value = "  +42.08%  "

# Row 40
$ws.Range("E40").Value = "  -4.97%  "

# Row 41
$ws.Range("E41").Value = "  -14.32%  "

# Row 42
$ws.Range("D42").Value = "2.72"
$ws.Range("E42").Value = "  -9.78%  "

# Row 43
$ws.Range("D43").Value = "0.354"
$ws.Range("E43").Value = "  -8.13%  "

# Row 44
$ws.Range("E44").Value = "  -7.62%  "

# Row 45
$ws.Range("D45").Value = "32.79"
$ws.Range("E45").Value = "  -10.81%  "

# Row 46
$ws.Range("E46").Value = "  -10.09%  "

# Row 47
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").Value = "2.69"
$ws.Range("E47").Value = "  -12.26%  "

# Row 48
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "3.16"
$ws.Range("E48").Value = "  -7.51%  "

# Row 49
$ws.Range("E49").Value = "  -8.00%  "

# Row 50
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  -0.43%  "

# Row 51
$ws.Range("D51").Value = "7.78"
$ws.Range("E51").Value = "  -9.76%  "
